# Update cryptocurrency price/volume data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.780.49'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.50%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.783.27'
$ws.Range("D3").Style = "Normal"

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.13%  '

$ws.Range("E6").Value = '  +0.00%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5122'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.91%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3847'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.11%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07832'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -7.95%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.087'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.23%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '40.69'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.73%  '

$ws.Range("E12").Value = '  +0.01%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.198'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.77%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.13'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.32%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.773.89'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.41%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.202'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.17%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.48'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.40%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001076'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.80%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06556'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.57%  '

$ws.Range("E20").Value = '  +0.03%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.02'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.22%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.910'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.97%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.826.59'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.47%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.01'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.79%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.233'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.82%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.90'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.30%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.21'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.99%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.986.37'
$ws.Range("D28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.362'
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.68'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.43%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1071'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.77%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.034'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.58%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.634'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.485'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.20%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07089'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.64%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '8.824'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.12%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02308'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.16%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2121'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.98%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.48'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.81%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.000'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.11%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6093'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.50%  '

$ws.Range("E42").Value = '  -0.02%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.152'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.55%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.319'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.80%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.12'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.24%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5882'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.23%  '

$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.707'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.06%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.91'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.40%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.198'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.44%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.897'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.60%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06844'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.89%  '
